# Revise the footer date "11th May 2011" -> "25th May 2011"
# The original single run containing "11th May " needs to become two runs:
#   "25"  and  "th May "
# We achieve the same visible result (and the same run-split) by locating
# the "11" at the start of that text and replacing it with "25", which
# naturally splits the run at the point of the edit.

$d = $word.ActiveDocument

foreach ($sec in $d.Sections) {
    $ftr = $sec.Footers.Item(1)  # wdHeaderFooterPrimary = 1
    if ($ftr.Exists) {
        $rng = $ftr.Range.Duplicate
        $rng.Find.ClearFormatting()
        $rng.Find.Text = "11th May "
        if ($rng.Find.Execute()) {
            # Select just the "11" portion (first two characters) and replace it.
            $numRange = $ftr.Range.Duplicate
            $numRange.Start = $rng.Start
            $numRange.End = $rng.Start + 2
            $numRange.Text = "25"
        }
    }
}
